# Bill of Materials update: resistor value changes + shift register part
# number change + new row for re-grouped resistors (shift register outputs
# 2.01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 40 (pushes the old row 40 and everything below it
# down by one, matching the new A1:E69 dimension).
$ws.Rows(40).Insert()

# Row 37 ("56" group): R68 and R71 moved in from the "39" group, R98 and
# R147 moved in from the "39" group; R70/R73/R100/R149 moved out to the
# new row 40.
$ws.Range("C37").Value = "R56, R58, R60, R62, R68, R71, R74, R75, R76, R77, R78, R79, R80, R81, R90, R91, R98, R147"
$ws.Range("D37").Value = "Resistor, Resistor, Resistor, Resistor, 56, 56, Resistor, Resistor, Resistor, Resistor, Resistor, Resistor, Resistor, Resistor, Resistor, Resistor, 56, 56"

# Row 38 ("39" group): R68, R71, R98, R147 moved out to the "56" group,
# leaving 14 designators.
$ws.Range("A38").Value = 14
$ws.Range("C38").Value = "R66, R82, R83, R84, R85, R86, R87, R88, R89, R95, R96, R97, R133, R134"
$ws.Range("D38").Value = "Resistor"

# Row 39 (R69, R72, R99, R148): resistor value changed from 150 to 200.
$ws.Range("B39").Value = 200
$ws.Range("D39").Value = 200

# New row 40: R70, R73, R100, R149 split out into their own 82-ohm group.
$ws.Range("A40").Value = 4
$ws.Range("B40").Value = 82
$ws.Range("C40").Value = "R70, R73, R100, R149"
$ws.Range("D40").Value = 82
$ws.Range("E40").Value = "RESC1608X06L"

# U13-U17 shift register part number change.
$ws.Range("B60").Value = "SN74HC595PWR"
